$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look like plain numbers remain stored as text
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.636.21'
$ws.Range('E2').Value = '  +1.18%  '

$ws.Range('D3').Value = '1.824.95'
$ws.Range('E3').Value = '  +1.85%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = '1.007'
$ws.Range('E5').Value = '  -0.01%  '

$ws.Range('D6').Value = '307.38'
$ws.Range('E6').Value = '  -0.04%  '

$ws.Range('D7').Value = '0.4634'
$ws.Range('E7').Value = '  +2.37%  '

$ws.Range('E8').Value = '  +0.30%  '

$ws.Range('D9').Value = '0.07130'
$ws.Range('E9').Value = '  +0.42%  '

$ws.Range('D10').Value = '0.9029'
$ws.Range('E10').Value = '  +2.09%  '

$ws.Range('D11').Value = '0.07769'
$ws.Range('E11').Value = '  -0.66%  '

$ws.Range('D12').Value = '19.37'
$ws.Range('E12').Value = '  -0.53%  '

$ws.Range('D13').Value = '1.804.40'
$ws.Range('E13').Value = '  +0.30%  '

$ws.Range('D14').Value = '5.264'
$ws.Range('E14').Value = '  -0.40%  '

$ws.Range('D15').Value = '6.343'
$ws.Range('E15').Value = '  +0.24%  '

$ws.Range('D16').Value = '87.80'
$ws.Range('E16').Value = '  +3.71%  '

$ws.Range('D17').Value = '1.008'

$ws.Range('D18').Value = '0.000008561'
$ws.Range('E18').Value = '  +0.26%  '

$ws.Range('E19').Value = '  -0.07%  '

$ws.Range('D20').Value = '26.667.44'
$ws.Range('E20').Value = '  +1.17%  '

$ws.Range('D21').Value = '14.17'
$ws.Range('E21').Value = '  -0.57%  '

$ws.Range('D22').Value = '5.019'
$ws.Range('E22').Value = '  +0.68%  '

$ws.Range('E23').Value = '  +0.33%  '

$ws.Range('D24').Value = '1.924'
$ws.Range('E24').Value = '  -2.60%  '

$ws.Range('D25').Value = '152.34'
$ws.Range('E25').Value = '  +0.16%  '

$ws.Range('D26').Value = '17.90'
$ws.Range('E26').Value = '  +0.15%  '

$ws.Range('E27').Value = '  -3.14%  '

$ws.Range('D28').Value = '113.83'
$ws.Range('E28').Value = '  +1.68%  '

$ws.Range('D29').Value = '4.830'
$ws.Range('E29').Value = '  -0.39%  '

$ws.Range('D30').Value = '0.08806'
$ws.Range('E30').Value = '  +1.53%  '

$ws.Range('D31').Value = '3.140'
$ws.Range('E31').Value = '  +2.60%  '

$ws.Range('D32').Value = '0.7321'
$ws.Range('E32').Value = '  +1.09%  '

$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').Value = '2.730'
$ws.Range('E33').Value = '  +0.97%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.438'
$ws.Range('E34').Value = '  -0.16%  '

$ws.Range('D35').Value = '1.138'
$ws.Range('E35').Value = '  +2.91%  '

$ws.Range('D36').Value = '1.075'
$ws.Range('E36').Value = '  +0.33%  '

$ws.Range('D37').Value = '0.01923'
$ws.Range('E37').Value = '  -0.31%  '

$ws.Range('D38').Value = '2.928'
$ws.Range('E38').Value = '  +1.84%  '

$ws.Range('D39').Value = '0.05127'
$ws.Range('E39').Value = '  +0.54%  '

$ws.Range('D40').Value = '6.874'
$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('D41').Value = '0.5058'
$ws.Range('E41').Value = '  -0.19%  '

$ws.Range('D42').Value = '0.1494'
$ws.Range('E42').Value = '  -1.30%  '

$ws.Range('D43').Value = '7.994'
$ws.Range('E43').Value = '  +0.15%  '

$ws.Range('D44').Value = '1.007'
$ws.Range('E44').Value = '  +0.00%  '

$ws.Range('D45').Value = '0.4660'
$ws.Range('E45').Value = '  +0.68%  '

$ws.Range('D46').Value = '9.981'
$ws.Range('E46').Value = '  +1.65%  '

$ws.Range('D47').Value = '98.41'
$ws.Range('E47').Value = '  -2.29%  '

$ws.Range('D48').Value = '1.559'
$ws.Range('E48').Value = '  -1.29%  '

$ws.Range('D49').Value = '0.05992'
$ws.Range('E49').Value = '  +0.41%  '

$ws.Range('D50').Value = '63.82'
$ws.Range('E50').Value = '  -0.47%  '

$ws.Range('D51').Value = '35.86'
$ws.Range('E51').Value = '  -0.65%  '
